$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header typos: "SDTree_D10" -> "DSTree_D10"
$ws.Range("F1").Value = "Fit time DSTree_D10"
$ws.Range("G1").Value = "Prediction time DSTree_D10"
$ws.Range("H1").Value = "Score DSTree_D10"

# Update numeric results for rows 2-8 (F:K)

# Row 2 - Iris
$ws.Range("F2").Value = 0.007318735122680664
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.9333333333333333
$ws.Range("I2").Value = 6.34464955329895
$ws.Range("J2").Value = 0.002009868621826172
$ws.Range("K2").Value = 0.9

# Row 3 - Wine
$ws.Range("F3").Value = 0.01031112670898438
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0.8888888888888888
$ws.Range("I3").Value = 0.007329702377319336
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0.9444444444444444

# Row 4 - Breast Cancer
$ws.Range("F4").Value = 0.03126907348632812
$ws.Range("G4").Value = 0.002786636352539062
$ws.Range("H4").Value = 0.9298245614035088
$ws.Range("I4").Value = 0.01590800285339355
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.9473684210526315

# Row 5 - Digits
$ws.Range("F5").Value = 0.1158242225646973
$ws.Range("G5").Value = 0.03172445297241211
$ws.Range("H5").Value = 0.7027777777777777
$ws.Range("I5").Value = 0.09540414810180664
$ws.Range("J5").Value = 0.01597309112548828
$ws.Range("K5").Value = 0.6722222222222223

# Row 6 - BankNote Authentication
$ws.Range("F6").Value = 0.02765989303588867
$ws.Range("G6").Value = 0.007353782653808594
$ws.Range("H6").Value = 0.92
$ws.Range("I6").Value = 0.02443718910217285
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0.9490909090909091

# Row 7 - Gas Drift
$ws.Range("F7").Value = 2.404195785522461
$ws.Range("G7").Value = 0.1112713813781738
$ws.Range("H7").Value = 0.8914450035945363
$ws.Range("I7").Value = 7.670363426208496
$ws.Range("J7").Value = 4.025373458862305
$ws.Range("K7").Value = 0.8914450035945363

# Row 8 - Shuttle
$ws.Range("F8").Value = 0.350999116897583
$ws.Range("G8").Value = 0.04653596878051758
$ws.Range("H8").Value = 0.9906034482758621
$ws.Range("I8").Value = 4.064779996871948
$ws.Range("J8").Value = 0.01574206352233887
$ws.Range("K8").Value = 0.9925862068965517

# Remove rows 9-13 (Adult, Hepatitis, Diabetes Progression, Boston Housing, California Housing)
$ws.Rows.Item(9).Resize(5).Delete()
